$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "life-dev/main"
$ws.Cells.Item(2, 3).Value = ""
$ws.Cells.Item(2, 4).Value = "scroll"
$ws.Cells.Item(2, 5).Value = ""
$ws.Cells.Item(2, 6).Value = "channel, page_url, scroll_rate, os_name"
$ws.Cells.Item(2, 7).Value = "Rround, https://life-dev.hectoinnovation.co.kr/main, 75, iOS"
$ws.Cells.Item(2, 8).Value = 4

# Row 3
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "life-dev/main"
$ws.Cells.Item(3, 3).Value = "상품"
$ws.Cells.Item(3, 4).Value = "click"
$ws.Cells.Item(3, 5).Value = "(스페인직수입)소르바스 포도씨유500ml 1P(스페인직수입)소르바스 포도씨유500ml 1P"
$ws.Cells.Item(3, 6).Value = "channel, page_url, click_text, module_id, module_order, prd_order, prd_code, prd_name, prd_brand, prd_price_origin, prd_price_final, prd_disc_rate, prd_is_ad, el_order, module_name, os_name"
$ws.Cells.Item(3, 7).Value = "Rround, https://life-dev.hectoinnovation.co.kr/main, (스페인직수입)소르바스 포도씨유500ml 1P(스페인직수입)소르바스 포도씨유500ml 1P, C-3, 13, 1, 3086, (스페인직수입)소르바스 포도씨유500ml 1P(스페인직수입)소르바스 포도씨유500ml 1P, 마이그스토어, 40,000원, 20,000원, 50%, F, 1, commerce-category-ranking, iOS"
$ws.Cells.Item(3, 8).Value = 16

# Row 4
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "ecommerce-dev/product/detail/3086"
$ws.Cells.Item(4, 3).Value = ""
$ws.Cells.Item(4, 4).Value = "pageview"
$ws.Cells.Item(4, 5).Value = ""
$ws.Cells.Item(4, 6).Value = "channel, page_url, prd_code, os_name"
$ws.Cells.Item(4, 7).Value = "Rround, https://ecommerce-dev.hectoinnovation.co.kr/product/detail/3086, 3086, iOS"
$ws.Cells.Item(4, 8).Value = 4

# Row 5
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "ecommerce-dev/product/detail/3086"
$ws.Cells.Item(5, 3).Value = ""
$ws.Cells.Item(5, 4).Value = "click"
$ws.Cells.Item(5, 5).Value = ""
$ws.Cells.Item(5, 6).Value = "channel, page_url, prd_code, os_name"
$ws.Cells.Item(5, 7).Value = "Rround, https://ecommerce-dev.hectoinnovation.co.kr/product/detail/3086, 3086, iOS"
$ws.Cells.Item(5, 8).Value = 4

# Row 6
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "life-dev/main"
$ws.Cells.Item(6, 3).Value = "뉴스"
$ws.Cells.Item(6, 4).Value = "click"
$ws.Cells.Item(6, 5).Value = "[OTT랭킹] '찰떡 캐스팅' 증명한 '광장'…K무비 침체 장기화"
$ws.Cells.Item(6, 6).Value = "channel, page_url, click_text, module_id, module_order, el_order, module_name, article_title, os_name"
$ws.Cells.Item(6, 7).Value = "Rround, https://life-dev.hectoinnovation.co.kr/main, [OTT랭킹] '찰떡 캐스팅' 증명한 '광장'…K무비 침체 장기화, D-1, 14, 1, news-card, [OTT랭킹] '찰떡 캐스팅' 증명한 '광장'…K무비 침체 장기화, iOS"
$ws.Cells.Item(6, 8).Value = 9

# Row 7
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "life-dev/news/detail/10736"
$ws.Cells.Item(7, 3).Value = ""
$ws.Cells.Item(7, 4).Value = "click"
$ws.Cells.Item(7, 5).Value = ""
$ws.Cells.Item(7, 6).Value = "channel, page_url, os_name"
$ws.Cells.Item(7, 7).Value = "Rround, https://life-dev.hectoinnovation.co.kr/news/detail/10736, iOS"
$ws.Cells.Item(7, 8).Value = 3

# Row 8
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "life-dev/main"
$ws.Cells.Item(8, 3).Value = "상품"
$ws.Cells.Item(8, 4).Value = "click"
$ws.Cells.Item(8, 5).Value = "(스페인직수입)소르바스 포도씨유500ml 1P(스페인직수입)소르바스 포도씨유500ml 1P"
$ws.Cells.Item(8, 6).Value = "channel, page_url, click_text, module_id, module_order, prd_order, prd_code, prd_name, prd_brand, prd_price_origin, prd_price_final, prd_disc_rate, prd_is_ad, el_order, module_name, os_name"
$ws.Cells.Item(8, 7).Value = "Rround, https://life-dev.hectoinnovation.co.kr/main, (스페인직수입)소르바스 포도씨유500ml 1P(스페인직수입)소르바스 포도씨유500ml 1P, C-3, 13, 1, 3086, (스페인직수입)소르바스 포도씨유500ml 1P(스페인직수입)소르바스 포도씨유500ml 1P, 마이그스토어, 40,000원, 20,000원, 50%, F, 1, commerce-category-ranking, iOS"
$ws.Cells.Item(8, 8).Value = 16

# Row 9
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "life-dev/main"
$ws.Cells.Item(9, 3).Value = "뉴스"
$ws.Cells.Item(9, 4).Value = "click"
$ws.Cells.Item(9, 5).Value = "[OTT랭킹] '찰떡 캐스팅' 증명한 '광장'…K무비 침체 장기화"
$ws.Cells.Item(9, 6).Value = "channel, page_url, click_text, module_id, module_order, el_order, module_name, article_title, os_name"
$ws.Cells.Item(9, 7).Value = "Rround, https://life-dev.hectoinnovation.co.kr/main, [OTT랭킹] '찰떡 캐스팅' 증명한 '광장'…K무비 침체 장기화, D-1, 14, 1, news-card, [OTT랭킹] '찰떡 캐스팅' 증명한 '광장'…K무비 침체 장기화, iOS"
$ws.Cells.Item(9, 8).Value = 9

# Row 10
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "life-dev/news/detail/10736"
$ws.Cells.Item(10, 3).Value = ""
$ws.Cells.Item(10, 4).Value = "click"
$ws.Cells.Item(10, 5).Value = ""
$ws.Cells.Item(10, 6).Value = "channel, page_url, os_name"
$ws.Cells.Item(10, 7).Value = "Rround, https://life-dev.hectoinnovation.co.kr/news/detail/10736, iOS"
$ws.Cells.Item(10, 8).Value = 3

# Row 11
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "life-dev/main"
$ws.Cells.Item(11, 3).Value = "상품"
$ws.Cells.Item(11, 4).Value = "click"
$ws.Cells.Item(11, 5).Value = "(스페인직수입)소르바스 포도씨유500ml 1P(스페인직수입)소르바스 포도씨유500ml 1P"
$ws.Cells.Item(11, 6).Value = "channel, page_url, click_text, module_id, module_order, prd_order, prd_code, prd_name, prd_brand, prd_price_origin, prd_price_final, prd_disc_rate, prd_is_ad, el_order, module_name, os_name"
$ws.Cells.Item(11, 7).Value = "Rround, https://life-dev.hectoinnovation.co.kr/main, (스페인직수입)소르바스 포도씨유500ml 1P(스페인직수입)소르바스 포도씨유500ml 1P, C-3, 13, 1, 3086, (스페인직수입)소르바스 포도씨유500ml 1P(스페인직수입)소르바스 포도씨유500ml 1P, 마이그스토어, 40,000원, 20,000원, 50%, F, 1, commerce-category-ranking, iOS"
$ws.Cells.Item(11, 8).Value = 16

# Row 12
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "ecommerce-dev/product/detail/3086"
$ws.Cells.Item(12, 3).Value = ""
$ws.Cells.Item(12, 4).Value = "pageview"
$ws.Cells.Item(12, 5).Value = ""
$ws.Cells.Item(12, 6).Value = "channel, page_url, prd_code, prd_price_origin, prd_price_final, prd_disc_rate, prd_review_cnt, prd_review_score, prd_tag, os_name"
$ws.Cells.Item(12, 7).Value = "Rround, https://ecommerce-dev.hectoinnovation.co.kr/product/detail/3086, 3086, 40,000원, 18,000원, 55%, 0, 0, #포도씨유___#올리브유___#소르바스___#압착오일___#엑스트라버진___#해바라기유___#카놀라유___#유기농, iOS"
$ws.Cells.Item(12, 8).Value = 10

# Row 13
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "ecommerce-dev/product/detail/3086"
$ws.Cells.Item(13, 3).Value = ""
$ws.Cells.Item(13, 4).Value = "click"
$ws.Cells.Item(13, 5).Value = ""
$ws.Cells.Item(13, 6).Value = "channel, page_url, tab_name, prd_code, prd_name, prd_price_origin, prd_price_final, prd_disc_rate, prd_review_cnt, prd_review_score, prd_tag, os_name"
$ws.Cells.Item(13, 7).Value = "Rround, https://ecommerce-dev.hectoinnovation.co.kr/product/detail/3086, 상품상세`n, 3086, (스페인직수입)소르바스 포도씨유500ml 1P(스페인직수입)소르바스 포도씨유500ml 1P, 40,000원, 20,000원, 55%, 0, 0, #포도씨유___#올리브유___#소르바스___#압착오일___#엑스트라버진___#해바라기유___#카놀라유___#유기농, iOS"
$ws.Cells.Item(13, 8).Value = 12

# Row 14
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "life-dev/main"
$ws.Cells.Item(14, 3).Value = ""
$ws.Cells.Item(14, 4).Value = "scroll"
$ws.Cells.Item(14, 5).Value = ""
$ws.Cells.Item(14, 6).Value = "channel, page_url, scroll_rate, os_name"
$ws.Cells.Item(14, 7).Value = "Rround, https://life-dev.hectoinnovation.co.kr/main, 75, iOS"
$ws.Cells.Item(14, 8).Value = 4

# Row 15
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "life-dev/main"
$ws.Cells.Item(15, 3).Value = "뉴스"
$ws.Cells.Item(15, 4).Value = "click"
$ws.Cells.Item(15, 5).Value = "K뮤지컬 통했다…'어쩌면 해피엔딩', 토니상 극본상·음악상 수상"
$ws.Cells.Item(15, 6).Value = "channel, page_url, click_text, module_id, module_order, el_order, module_name, article_title, os_name"
$ws.Cells.Item(15, 7).Value = "Rround, https://life-dev.hectoinnovation.co.kr/main, K뮤지컬 통했다…'어쩌면 해피엔딩', 토니상 극본상·음악상 수상, D-1, 19, 1, news-card, K뮤지컬 통했다…'어쩌면 해피엔딩', 토니상 극본상·음악상 수상, iOS"
$ws.Cells.Item(15, 8).Value = 9
